$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new year column (2022), mirroring the formatting of column R
# (the last existing year column) for the header + two data rows.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 30
$ws.Range("S6").Value = 11928.6

# Update the active selection as in the edited workbook
$ws.Range("T3").Select()
